$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = '":null,"thumbhttps://m.media-amazon.com/images/I/51umk+HPQbL._AC_SR38,50_.jpg'
$ws.Range("B1").Value = "iphone 7 negro"
